$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.670.89"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").Value = "'2.058.87"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'243.56"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "'0.668"
$ws.Range("E6").Value = "  +2.12%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'54.41"
$ws.Range("E8").Value = "  -6.83%  "

$ws.Range("D9").Value = "'59.06"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").Value = "'0.365"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").Value = "'0.0751"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("E12").Value = "  -2.99%  "

$ws.Range("D13").Value = "'0.937"
$ws.Range("E13").Value = "  +6.21%  "

$ws.Range("D14").Value = "'14.76"
$ws.Range("E14").Value = "  -3.53%  "

$ws.Range("D15").Value = "'2.358.71"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("E16").Value = "  -2.91%  "

$ws.Range("D17").Value = "'2.066.50"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("D18").Value = "'36.567.82"
$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").Value = "'17.03"
$ws.Range("E19").Value = "  -6.49%  "

$ws.Range("D20").Value = "'72.01"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("D21").Value = "'0.0₃0863"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").Value = "'237.96"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("D23").Value = "'5.26"
$ws.Range("E23").Value = "  -1.78%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("D28").Value = "'164.44"
$ws.Range("E28").Value = "  -1.99%  "

$ws.Range("D29").Value = "'20.14"
$ws.Range("E29").Value = "  +1.32%  "

$ws.Range("E30").Value = "  -1.17%  "

$ws.Range("E31").Value = "  -8.31%  "

$ws.Range("E32").Value = "  +7.75%  "

$ws.Range("D33").Value = "'4.51"
$ws.Range("E33").Value = "  -5.21%  "

$ws.Range("D34").Value = "'0.0597"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("D38").Value = "'0.0832"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("D40").Value = "'4.94"
$ws.Range("E40").Value = "  -5.44%  "

$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  -6.71%  "

$ws.Range("D42").Value = "'0.0216"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").Value = "'1.11"
$ws.Range("E43").Value = "  -2.87%  "

$ws.Range("D44").Value = "'94.44"
$ws.Range("E44").Value = "  -2.56%  "

$ws.Range("D45").Value = "'0.0911"
$ws.Range("E45").Value = "  -3.69%  "

$ws.Range("D46").Value = "'1.409.80"
$ws.Range("E46").Value = "  +9.21%  "

$ws.Range("E47").Value = "  +14.10%  "

$ws.Range("D48").Value = "'16.02"
$ws.Range("E48").Value = "  -4.83%  "

$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  +2.24%  "

$ws.Range("E50").Value = "  -3.61%  "

$ws.Range("D51").Value = "'2.246.69"
$ws.Range("E51").Value = "  +1.29%  "

